$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 2: Mercury Rising | Quicksilver
$ws.Range("H2").Value = 1430.1818
$ws.Range("I2").Value = 1591.5
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1591.5
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -1478.5
$ws.Range("N2").Value = -1226

# ALC row 19: Unbreak My Heart | Roof Tile
$ws.Range("H19").Value = 3554.6
$ws.Range("I19").Value = 294.75
$ws.Range("J19").Value = 4740
$ws.Range("K19").Value = 294.75
$ws.Range("L19").Value = 4740
$ws.Range("M19").Value = -119.75
$ws.Range("N19").Value = -5090

# ALC row 32: Automata for the People | Crab Oil
$ws.Range("H32").Value = 2099.889
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 2112.375
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 2112.375
$ws.Range("M32").Value = -1674
$ws.Range("N32").Value = -2764.375

# ALC row 47: Open Your Grimoire to Page 42 | Embossed Book of Silver
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

# ALC row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 17860284
$ws.Range("I138").Value = 1503.4
$ws.Range("J138").Value = 38466570
$ws.Range("K138").Value = 4510.200000000001
$ws.Range("L138").Value = 115399710
$ws.Range("M138").Value = 629.7999999999993
$ws.Range("N138").Value = -115409990

$ws = $wb.Worksheets.Item("ARM")
# ARM row 5: The Alloyed Truth | Bronze Rivets
$ws.Range("H5").Value = 2390
$ws.Range("I5").Value = 2650
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 2650
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = -2538

# ARM row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 4519.6313
$ws.Range("I32").Value = 2959.6365
$ws.Range("J32").Value = 9799.615
$ws.Range("K32").Value = 2959.6365
$ws.Range("L32").Value = 9799.615
$ws.Range("M32").Value = -2672.6365

# ARM row 52: Distill and Know that I'm Right | Mythril Alembic
$ws.Range("H52").Value = 39779.09
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 39779.09
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 39779.09
$ws.Range("N52").Value = -40415.09

# ARM row 92: Mail It In | High Steel Scale Mail of Fending
$ws.Range("H92").Value = 31071.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 31071.5
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 31071.5
$ws.Range("N92").Value = -36063.5

$ws = $wb.Worksheets.Item("BSM")
# BSM row 4: Mending Fences | Bronze Rivets
$ws.Range("H4").Value = 2390
$ws.Range("I4").Value = 2650
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 2650
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = -2535

# BSM row 82: Spirituality Inspector | Titanium Lump Hammer
$ws.Range("H82").Value = 7230
$ws.Range("I82").Value = 2651
$ws.Range("J82").Value = 39283
$ws.Range("K82").Value = 2651
$ws.Range("L82").Value = 39283
$ws.Range("M82").Value = -2268

# BSM row 85: The Clamor for Hammers (L) | Titanium Lump Hammer
$ws.Range("H85").Value = 7230
$ws.Range("I85").Value = 2651
$ws.Range("J85").Value = 39283
$ws.Range("K85").Value = 2651
$ws.Range("L85").Value = 39283
$ws.Range("M85").Value = -1325

# BSM row 92: Have Blade, Will Travel | High Steel Katzbalger
$ws.Range("H92").Value = 20401
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 20401
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 20401
$ws.Range("N92").Value = -25393

$ws = $wb.Worksheets.Item("CRP")
# CRP row 68: Do You Even String Bow | Holy Cedar Composite Bow
$ws.Range("H68").Value = 48000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 48000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 48000
$ws.Range("N68").Value = -49498

# CRP row 71: Win One Bow, Get Three Free (L) | Holy Cedar Composite Bow
$ws.Range("H71").Value = 48000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 48000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 144000
$ws.Range("N71").Value = -151488

$ws = $wb.Worksheets.Item("CUL")
# CUL row 4: In Hot Water | Boiled Egg
$ws.Range("H4").Value = 2013.0769
$ws.Range("I4").Value = 92.57143000000001
$ws.Range("J4").Value = 4253.6665
$ws.Range("K4").Value = 277.71429
$ws.Range("L4").Value = 12760.9995
$ws.Range("M4").Value = -165.71429
$ws.Range("N4").Value = -12984.9995

# CUL row 63: The Next to Last Supper | Stuffed Cabbage Rolls
$ws.Range("H63").Value = 2502.4
$ws.Range("I63").Value = 1004.8
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 3014.4
$ws.Range("L63").Value = 12000
$ws.Range("M63").Value = -2265.4
$ws.Range("N63").Value = -13498

# CUL row 64: The Aroma of Faith | Baked Onion Soup
$ws.Range("H64").Value = 1690
$ws.Range("I64").Value = 1238
$ws.Range("J64").Value = 2012.8572
$ws.Range("K64").Value = 3714
$ws.Range("L64").Value = 6038.571599999999
$ws.Range("M64").Value = -3444
$ws.Range("N64").Value = -6578.571599999999

# CUL row 66: Nostalgia through the Stomach (L) | Stuffed Cabbage Rolls
$ws.Range("H66").Value = 2502.4
$ws.Range("I66").Value = 1004.8
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 9043.199999999999
$ws.Range("L66").Value = 36000
$ws.Range("M66").Value = -5299.199999999999
$ws.Range("N66").Value = -43488

# CUL row 67: Soup's On (L) | Baked Onion Soup
$ws.Range("H67").Value = 1690
$ws.Range("I67").Value = 1238
$ws.Range("J67").Value = 2012.8572
$ws.Range("K67").Value = 3714
$ws.Range("L67").Value = 6038.571599999999
$ws.Range("M67").Value = -2778
$ws.Range("N67").Value = -7910.571599999999

# CUL row 81: It Goes Down Smoothly | Frozen Spirits
$ws.Range("H81").Value = 36000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 36000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 108000
$ws.Range("N81").Value = -110246

# CUL row 84: Quenching the Flame (L) | Frozen Spirits
$ws.Range("H84").Value = 36000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 36000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 324000
$ws.Range("N84").Value = -335232

# CUL row 106: Herky Jerky | Jerked Jhammel
$ws.Range("H106").Value = 4050
$ws.Range("I106").Value = 1000
$ws.Range("J106").Value = 4485.7144
$ws.Range("K106").Value = 3000
$ws.Range("L106").Value = 13457.1432
$ws.Range("M106").Value = -2054
$ws.Range("N106").Value = -15349.1432

# CUL row 134: Don't Knock It Till You've Tried It | Mezcal-marinated Swampmonk
$ws.Range("H134").Value = 1512.5294
$ws.Range("I134").Value = 1086.6666
$ws.Range("J134").Value = 1991.625
$ws.Range("K134").Value = 3259.9998
$ws.Range("L134").Value = 5974.875
$ws.Range("M134").Value = 1810.0002
$ws.Range("N134").Value = -16114.875

$ws = $wb.Worksheets.Item("GSM")
# GSM row 2: Copper and Robbers | Copper Ingot
$ws.Range("H2").Value = 64.84614999999999
$ws.Range("I2").Value = 67.333336
$ws.Range("J2").Value = 59.25
$ws.Range("K2").Value = 67.333336
$ws.Range("L2").Value = 59.25
$ws.Range("M2").Value = 45.666664

# GSM row 18: Gorgeous Gorget | Brass Gorget
$ws.Range("H18").Value = 70006
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 70006
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 70006
$ws.Range("N18").Value = -70592

# GSM row 22: Bad to the Bone | Brass Circlet (Sunstone)
$ws.Range("H22").Value = 70009
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 70009
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 70009
$ws.Range("N22").Value = -71067
$ws.Range("M22").ClearContents()

# GSM row 46: Burning the Midnight Oil | Fire Brand
$ws.Range("H46").Value = 2000
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1844
$ws.Range("N46").ClearContents()

# GSM row 80: Needs More Prayerbell | Hardsilver Ingot
$ws.Range("H80").Value = 4542.7144
$ws.Range("I80").Value = 4542.7144
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4542.7144
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3544.7144
$ws.Range("N80").ClearContents()

# GSM row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws.Range("H83").Value = 4542.7144
$ws.Range("I83").Value = 4542.7144
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 22713.572
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -17721.572
$ws.Range("N83").ClearContents()

# GSM row 97: If I'd a Koppranickel for Every Time... | Koppranickel Ingot
$ws.Range("H97").Value = 1987.25
$ws.Range("I97").Value = 842.5714
$ws.Range("J97").Value = 10000
$ws.Range("K97").Value = 842.5714
$ws.Range("L97").Value = 10000
$ws.Range("M97").Value = -346.5714

# GSM row 107: Whetstones for the Workers | Hard Mudstone Whetstone
$ws.Range("H107").Value = 952.7143
$ws.Range("I107").Value = 224.28572
$ws.Range("J107").Value = 1681.1428
$ws.Range("K107").Value = 224.28572
$ws.Range("L107").Value = 1681.1428
$ws.Range("M107").Value = 1695.71428
$ws.Range("N107").Value = -5521.1428

$ws = $wb.Worksheets.Item("LTW")
# LTW row 16: Saddle Sore | Hard Leather
$ws.Range("H16").Value = 3116.5417
$ws.Range("I16").Value = 1599.8334
$ws.Range("J16").Value = 7666.6665
$ws.Range("K16").Value = 1599.8334
$ws.Range("L16").Value = 7666.6665
$ws.Range("M16").Value = -1429.8334
$ws.Range("N16").Value = -8006.6665

# LTW row 21: Heads Up | Hard Leather Skullcap
$ws.Range("H21").Value = 46782.332
$ws.Range("I21").Value = 333
$ws.Range("J21").Value = 70007
$ws.Range("K21").Value = 333
$ws.Range("L21").Value = 70007
$ws.Range("M21").Value = -159
$ws.Range("N21").Value = -70355

$ws = $wb.Worksheets.Item("WVR")
# WVR row 61: Bundle Up, It's Odd out There | Woolen Deerstalker
$ws.Range("H61").Value = 62542.75
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 62542.75
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 62542.75
$ws.Range("N61").Value = -63126.75
